# Update the Training Dashboard with the 04-Nov-2025 progress snapshot.
# For each data row (3..25):
#   - PERIOD TO EXPIRE (col H) decreases by 1 day
#   - LAST UPDATE (col I) text changes from "03-Nov-2025" to "04-Nov-2025"
#
# The LAST UPDATE column stores the date as a literal text string (not a
# real Excel date value), so we must avoid Excel's automatic "looks like a
# date" type coercion when writing the new text - otherwise the cell would
# turn into a date serial number and pick up a brand-new date number
# format/style. We work around this by building the text in an unused
# scratch cell via a text formula (="04-Nov-2025"), which always yields a
# plain text result, and then copying only the *value* of that cell onto
# the target cell - the target cell's existing style is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$oldDate = "03-Nov-2025"
$newDate = "04-Nov-2025"

# Far-away scratch cell, outside the used range, used to mint a text value
# without triggering date auto-detection.
$scratch = $ws.Cells.Item(500, 500)
$scratch.Formula = "=""" + $newDate + """"
$scratch.Copy()

for ($r = 3; $r -le 25; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $hCell.Value() - 1

    $iCell = $ws.Cells.Item($r, 9)
    $iCell.PasteSpecial(-4163)   # xlPasteValues - copies only the text, keeps I-cell's style
}

$scratch.Clear()
$excel.CutCopyMode = 0
